$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 6283
$ws.Range("K3").Value = 6485
$ws.Range("I4").Value = 1486
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 1357
$ws.Range("K5").Value = 462
$ws.Range("K6").Value = 7150
$ws.Range("I7").Value = 20784
$ws.Range("J7").Value = 23206
$ws.Range("K7").Value = 21737

$ws = $wb.Worksheets.Item("Grant Park")
$ws.Range("K5").Value = 12
$ws.Range("K6").Value = 23

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 389
$ws.Range("K3").Value = 435
$ws.Range("K4").Value = 80
$ws.Range("K5").Value = 43
$ws.Range("K6").Value = 480
$ws.Range("K7").Value = 1427

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 245
$ws.Range("K5").Value = 25
$ws.Range("K6").Value = 290
$ws.Range("K7").Value = 946

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K6").Value = 220
$ws.Range("K7").Value = 738

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 92
$ws.Range("K7").Value = 358

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 189
$ws.Range("K7").Value = 639
$ws.Range("K8").Value = 1427
$ws.Range("K9").Value = 97
$ws.Range("K11").Value = 402
$ws.Range("K19").Value = 637
$ws.Range("K20").Value = 518
$ws.Range("K30").Value = 81
$ws.Range("K31").Value = 241
$ws.Range("K33").Value = 946
$ws.Range("K35").Value = 34
$ws.Range("K37").Value = 738
$ws.Range("K38").Value = 23
$ws.Range("K42").Value = 802
$ws.Range("K46").Value = 44
$ws.Range("K47").Value = 149
$ws.Range("K48").Value = 272
$ws.Range("K52").Value = 570
$ws.Range("K54").Value = 429
$ws.Range("K60").Value = 127
$ws.Range("I63").Value = 186
$ws.Range("J63").Value = 103
$ws.Range("K63").Value = 63
$ws.Range("K67").Value = 849
$ws.Range("K73").Value = 196
$ws.Range("K76").Value = 298
$ws.Range("K77").Value = 150
$ws.Range("K78").Value = 245
$ws.Range("K79").Value = 548
$ws.Range("K84").Value = 174
$ws.Range("K85").Value = 1012
$ws.Range("K88").Value = 231
$ws.Range("K90").Value = 203
$ws.Range("K94").Value = 291
$ws.Range("K98").Value = 108
$ws.Range("K99").Value = 358
$ws.Range("K100").Value = 41
$ws.Range("I101").Value = 20784
$ws.Range("J101").Value = 23206
$ws.Range("K101").Value = 21737

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 78
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 241

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 308
$ws.Range("K4").Value = 48
$ws.Range("K7").Value = 849

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 68
$ws.Range("K7").Value = 174

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 105
$ws.Range("K7").Value = 429

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K2").Value = 41
$ws.Range("K7").Value = 272

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K5").Value = 19
$ws.Range("K6").Value = 207
$ws.Range("K7").Value = 637

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K6").Value = 154
$ws.Range("K7").Value = 298

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K6").Value = 300
$ws.Range("K7").Value = 802

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 73
$ws.Range("K3").Value = 60
$ws.Range("K7").Value = 245

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("K2").Value = 19
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K3").Value = 177
$ws.Range("K7").Value = 548

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 169
$ws.Range("K7").Value = 518

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 41

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 209
$ws.Range("K6").Value = 176
$ws.Range("K7").Value = 639

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 130
$ws.Range("K7").Value = 291

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K3").Value = 41
$ws.Range("K7").Value = 149

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 63
$ws.Range("K7").Value = 108

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K3").Value = 103
$ws.Range("K7").Value = 402

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("K2").Value = 3
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 97

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 65
$ws.Range("K7").Value = 196

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K3").Value = 51
$ws.Range("K7").Value = 189

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K2").Value = 59
$ws.Range("K3").Value = 69
$ws.Range("K7").Value = 231

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K2").Value = 75
$ws.Range("K3").Value = 60
$ws.Range("K7").Value = 203

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K2").Value = 42
$ws.Range("K7").Value = 127

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 330
$ws.Range("K3").Value = 347
$ws.Range("K7").Value = 1012

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 59
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 150

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K4").Value = 32
$ws.Range("K6").Value = 204
$ws.Range("K7").Value = 570
